# Apply updated TPM-derived values (per commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Target cluster -> Inflammatory-Mac
$row2 = New-Object "object[,]" 1,17
$row2[0,0] = "Inflammatory-Mac"
$row2[0,1] = 3
$row2[0,2] = 1
$row2[0,3] = 1.808798
$row2[0,4] = 5.426394
$row2[0,5] = 0.36666838522954
$row2[0,6] = 0.4612006375427627
$row2[0,7] = 3
$row2[0,8] = 1
$row2[0,9] = 0.4260053333333333
$row2[0,10] = 1.278016
$row2[0,11] = 0.02405532912416773
$row2[0,12] = 0.02531756756689831
$row2[0,13] = 0.7705575949226667
$row2[0,14] = 6.935018354304001
$row2[0,15] = 0.008820328686123705
$row2[0,16] = 0.01167647830288547
$ws.Range("D2:T2").Value = $row2

# Row 3: Target cluster -> ECs
$row3 = New-Object "object[,]" 1,17
$row3[0,0] = "ECs"
$row3[0,1] = 3
$row3[0,2] = 1
$row3[0,3] = 1.808798
$row3[0,4] = 5.426394
$row3[0,5] = 0.36666838522954
$row3[0,6] = 0.4612006375427627
$row3[0,7] = 3
$row3[0,8] = 1
$row3[0,9] = 1.041192666666667
$row3[0,10] = 3.123578
$row3[0,11] = 0.05879323641880037
$row3[0,12] = 0.06187825274916518
$row3[0,13] = 1.883307213081334
$row3[0,14] = 16.949764917732
$row3[0,15] = 0.02155762106010011
$row3[0,16] = 0.02853828961794719
$ws.Range("D3:T3").Value = $row3

# Row 4: Target cluster -> Resolving-Mac
$row4 = New-Object "object[,]" 1,17
$row4[0,0] = "Resolving-Mac"
$row4[0,1] = 3
$row4[0,2] = 1
$row4[0,3] = 1.808798
$row4[0,4] = 5.426394
$row4[0,5] = 0.36666838522954
$row4[0,6] = 0.4612006375427627
$row4[0,7] = 3
$row4[0,8] = 1
$row4[0,9] = 5.850740666666667
$row4[0,10] = 17.552222
$row4[0,11] = 0.3303749538898241
$row4[0,12] = 0.3477104875323931
$row4[0,13] = 10.58280801638533
$row4[0,14] = 95.245272147468
$row4[0,15] = 0.1211380508630656
$row4[0,16] = 0.1603642985302446
$ws.Range("D4:T4").Value = $row4

# Row 5: Target cluster -> FAPs
$row5 = New-Object "object[,]" 1,17
$row5[0,0] = "FAPs"
$row5[0,1] = 3
$row5[0,2] = 1
$row5[0,3] = 1.808798
$row5[0,4] = 5.426394
$row5[0,5] = 0.36666838522954
$row5[0,6] = 0.4612006375427627
$row5[0,7] = 2
$row5[0,8] = 1
$row5[0,9] = 2.648771
$row5[0,10] = 5.297542
$row5[0,11] = 0.1495686865725097
$row5[0,12] = 0.1049445996947469
$row5[0,13] = 4.791091687258
$row5[0,14] = 28.746550123548
$row5[0,15] = 0.05484210878644533
$row5[0,16] = 0.04840051628588728
$ws.Range("D5:T5").Value = $row5

# Row 6: Target cluster -> MuSCs
$row6 = New-Object "object[,]" 1,17
$row6[0,0] = "MuSCs"
$row6[0,1] = 3
$row6[0,2] = 1
$row6[0,3] = 1.808798
$row6[0,4] = 5.426394
$row6[0,5] = 0.36666838522954
$row6[0,6] = 0.4612006375427627
$row6[0,7] = 3
$row6[0,8] = 1
$row6[0,9] = 7.742685666666667
$row6[0,10] = 23.228057
$row6[0,11] = 0.4372077939946981
$row6[0,12] = 0.4601490924567965
$row6[0,13] = 14.00495434849533
$row6[0,14] = 126.044589136458
$row6[0,15] = 0.1603102758338053
$row6[0,16] = 0.2122210548057982
$ws.Range("D6:T6").Value = $row6

# Row 7: Target cluster -> Inflammatory-Mac
$row7 = New-Object "object[,]" 1,17
$row7[0,0] = "Inflammatory-Mac"
$row7[0,1] = 2
$row7[0,2] = 1
$row7[0,3] = 3.033388
$row7[0,4] = 6.066776
$row7[0,5] = 0.6149097244328354
$row7[0,6] = 0.5156280504196953
$row7[0,7] = 3
$row7[0,8] = 1
$row7[0,9] = 0.4260053333333333
$row7[0,10] = 1.278016
$row7[0,11] = 0.02405532912416773
$row7[0,12] = 0.02531756756689831
$row7[0,13] = 1.292239466069333
$row7[0,14] = 7.753436796416
$row7[0,15] = 0.01479185580288314
$row7[0,16] = 0.01305444800588868
$ws.Range("D7:T7").Value = $row7

# Row 8: Target cluster -> ECs
$row8 = New-Object "object[,]" 1,17
$row8[0,0] = "ECs"
$row8[0,1] = 2
$row8[0,2] = 1
$row8[0,3] = 3.033388
$row8[0,4] = 6.066776
$row8[0,5] = 0.6149097244328354
$row8[0,6] = 0.5156280504196953
$row8[0,7] = 3
$row8[0,8] = 1
$row8[0,9] = 1.041192666666667
$row8[0,10] = 3.123578
$row8[0,11] = 0.05879323641880037
$row8[0,12] = 0.06187825274916518
$row8[0,13] = 3.158341340754667
$row8[0,14] = 18.950048044528
$row8[0,15] = 0.03615253280479908
$row8[0,16] = 0.0319061628284292
$ws.Range("D8:T8").Value = $row8

# Row 9: Target cluster -> Resolving-Mac
$row9 = New-Object "object[,]" 1,17
$row9[0,0] = "Resolving-Mac"
$row9[0,1] = 2
$row9[0,2] = 1
$row9[0,3] = 3.033388
$row9[0,4] = 6.066776
$row9[0,5] = 0.6149097244328354
$row9[0,6] = 0.5156280504196953
$row9[0,7] = 3
$row9[0,8] = 1
$row9[0,9] = 5.850740666666667
$row9[0,10] = 17.552222
$row9[0,11] = 0.3303749538898241
$row9[0,12] = 0.3477104875323931
$row9[0,13] = 17.74756652937867
$row9[0,14] = 106.485399176272
$row9[0,15] = 0.2031507718559025
$row9[0,16] = 0.1792892807968096
$ws.Range("D9:T9").Value = $row9

# Row 10: Target cluster -> FAPs
$row10 = New-Object "object[,]" 1,17
$row10[0,0] = "FAPs"
$row10[0,1] = 2
$row10[0,2] = 1
$row10[0,3] = 3.033388
$row10[0,4] = 6.066776
$row10[0,5] = 0.6149097244328354
$row10[0,6] = 0.5156280504196953
$row10[0,7] = 2
$row10[0,8] = 1
$row10[0,9] = 2.648771
$row10[0,10] = 5.297542
$row10[0,11] = 0.1495686865725097
$row10[0,12] = 0.1049445996947469
$row10[0,13] = 8.034750166147999
$row10[0,14] = 32.139000664592
$row10[0,15] = 0.09197123984408309
$row10[0,16] = 0.05411237934267767
$ws.Range("D10:T10").Value = $row10

# Row 11: Target cluster -> MuSCs
$row11 = New-Object "object[,]" 1,17
$row11[0,0] = "MuSCs"
$row11[0,1] = 2
$row11[0,2] = 1
$row11[0,3] = 3.033388
$row11[0,4] = 6.066776
$row11[0,5] = 0.6149097244328354
$row11[0,6] = 0.5156280504196953
$row11[0,7] = 3
$row11[0,8] = 1
$row11[0,9] = 7.742685666666667
$row11[0,10] = 23.228057
$row11[0,11] = 0.4372077939946981
$row11[0,12] = 0.4601490924567965
$row11[0,13] = 23.48656978903866
$row11[0,14] = 140.919418734232
$row11[0,15] = 0.2688433241251676
$row11[0,16] = 0.2372657794458901
$ws.Range("D11:T11").Value = $row11

# Row 12: Target cluster -> Inflammatory-Mac
$row12 = New-Object "object[,]" 1,17
$row12[0,0] = "Inflammatory-Mac"
$row12[0,1] = 1
$row12[0,2] = 0.3333333333333333
$row12[0,3] = 0.09087633333333334
$row12[0,4] = 0.272629
$row12[0,5] = 0.01842189033762463
$row12[0,6] = 0.02317131203754203
$row12[0,7] = 3
$row12[0,8] = 1
$row12[0,9] = 0.4260053333333333
$row12[0,10] = 1.278016
$row12[0,11] = 0.02405532912416773
$row12[0,12] = 0.02531756756689831
$row12[0,13] = 0.03871380267377778
$row12[0,14] = 0.348424224064
$row12[0,15] = 0.0004431446351608858
$row12[0,16] = 0.0005866412581241546
$ws.Range("D12:T12").Value = $row12

# Row 13: Target cluster -> ECs
$row13 = New-Object "object[,]" 1,17
$row13[0,0] = "ECs"
$row13[0,1] = 1
$row13[0,2] = 0.3333333333333333
$row13[0,3] = 0.09087633333333334
$row13[0,4] = 0.272629
$row13[0,5] = 0.01842189033762463
$row13[0,6] = 0.02317131203754203
$row13[0,7] = 3
$row13[0,8] = 1
$row13[0,9] = 1.041192666666667
$row13[0,10] = 3.123578
$row13[0,11] = 0.05879323641880037
$row13[0,12] = 0.06187825274916518
$row13[0,13] = 0.09461977184022223
$row13[0,14] = 0.851577946562
$row13[0,15] = 0.001083082553901179
$row13[0,16] = 0.0014338003027888
$ws.Range("D13:T13").Value = $row13

# Row 14: Target cluster -> Resolving-Mac
$row14 = New-Object "object[,]" 1,17
$row14[0,0] = "Resolving-Mac"
$row14[0,1] = 1
$row14[0,2] = 0.3333333333333333
$row14[0,3] = 0.09087633333333334
$row14[0,4] = 0.272629
$row14[0,5] = 0.01842189033762463
$row14[0,6] = 0.02317131203754203
$row14[0,7] = 3
$row14[0,8] = 1
$row14[0,9] = 5.850740666666667
$row14[0,10] = 17.552222
$row14[0,11] = 0.3303749538898241
$row14[0,12] = 0.3477104875323931
$row14[0,13] = 0.531693859070889
$row14[0,14] = 4.785244731638
$row14[0,15] = 0.006086131170856134
$row14[0,16] = 0.00805690820533895
$ws.Range("D14:T14").Value = $row14

# Row 15: Target cluster -> FAPs
$row15 = New-Object "object[,]" 1,17
$row15[0,0] = "FAPs"
$row15[0,1] = 1
$row15[0,2] = 0.3333333333333333
$row15[0,3] = 0.09087633333333334
$row15[0,4] = 0.272629
$row15[0,5] = 0.01842189033762463
$row15[0,6] = 0.02317131203754203
$row15[0,7] = 2
$row15[0,8] = 1
$row15[0,9] = 2.648771
$row15[0,10] = 5.297542
$row15[0,11] = 0.1495686865725097
$row15[0,12] = 0.1049445996947469
$row15[0,13] = 0.2407105963196667
$row15[0,14] = 1.444263577918
$row15[0,15] = 0.002755337941981324
$row15[0,16] = 0.002431704066181918
$ws.Range("D15:T15").Value = $row15

# Row 16: Target cluster -> MuSCs
$row16 = New-Object "object[,]" 1,17
$row16[0,0] = "MuSCs"
$row16[0,1] = 1
$row16[0,2] = 0.3333333333333333
$row16[0,3] = 0.09087633333333334
$row16[0,4] = 0.272629
$row16[0,5] = 0.01842189033762463
$row16[0,6] = 0.02317131203754203
$row16[0,7] = 3
$row16[0,8] = 1
$row16[0,9] = 7.742685666666667
$row16[0,10] = 23.228057
$row16[0,11] = 0.4372077939946981
$row16[0,12] = 0.4601490924567965
$row16[0,13] = 0.7036268835392222
$row16[0,14] = 6.332641951853001
$row16[0,15] = 0.008054194035725108
$row16[0,16] = 0.01066225820510821
$ws.Range("D16:T16").Value = $row16
